$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "ferreira.victor@mrv.com.br"
$ws.Range("B5").Value = "Ferramenta x"
$ws.Range("C5").Value = "Ferramenta de Planejamento"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "Ferrmanta 1"
